# Implement min up/down time for thermal generators
# - Insert two new columns (J:K) for "MinUpTime" / "MinDownTime" ([h]) right
#   after the existing RampDw column, shifting Qmax..long to the right.
# - Populate the new header cells and units row.
# - Populate MinUpTime/MinDownTime values for the OCGT rows that have them.
# - Update the named ranges / filter database that referenced the old
#   right-hand edge of the table (column W) so they now reach column Y.
# - Restore the active-cell selection on the frozen pane.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at J:K (existing Qmax..long data shifts right).
$ws.Range("J:K").Insert()

# Units row (row 6) for the two new columns -> "[h]"
$ws.Range("J6").Value = "[h]"
$ws.Range("K6").Value = "[h]"

# Header row (row 3) for the two new columns
$ws.Range("J3").Value = "MinUpTime"
$ws.Range("K3").Value = "MinDownTime"

# Data for the rows that specify Min Up/Down times
$ws.Range("J21").Value = 3
$ws.Range("K21").Value = 2

$ws.Range("J22").Value = 3

$ws.Range("K23").Value = 4

$ws.Range("J25").Value = 3
$ws.Range("K25").Value = 2

# Update named ranges that pointed at the old table extent ($W$26 -> $Y$26)
foreach ($n in $wb.Names) {
    if ($n.Name -eq "thermalgen") {
        $n.RefersTo = "='Power ThermalGen'!`$B`$7:`$Y`$26"
    }
    if ($n.Name -eq "Power ThermalGen!_FilterDatabase") {
        $n.RefersTo = "='Power ThermalGen'!`$B`$7:`$Y`$26"
    }
}

# Restore the selection on the frozen (bottom-left) pane
$ws.Range("J12").Select()

Write-Output "done"
